$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1) refreshed to a later pull time
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 10:22"

# Austria (row 20) - updated totals
$ws.Range("B20").Value = 13818
$ws.Range("C20").Value = 12
$ws.Range("E20").Value = 6877

# Australia (row 32) - updated totals
$ws.Range("B32").Value = 6313
$ws.Range("C32").Value = 10
$ws.Range("D32").Value = 3338
$ws.Range("E32").Value = 2916
$ws.Range("F32").Value = 81

# Malasia / Filipinas swap ranking places (row 37 / row 38)
$ws.Range("A37").Value = "Filipinas"
$ws.Range("B37").Value = 4648
$ws.Range("C37").Value = 220
$ws.Range("D37").Value = 197
$ws.Range("E37").Value = 4154
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 50
$ws.Range("H37").Value = 297

$ws.Range("A38").Value = "Malasia"
$ws.Range("B38").Value = 4530
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 1995
$ws.Range("E38").Value = 2462
$ws.Range("F38").Value = 72
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 73

# Bielorrusia moves ahead of Tailandia and Singapur (rows 51-53)
$ws.Range("A51").Value = "Bielorrusia"
$ws.Range("B51").Value = 2578
$ws.Range("C51").Value = 352
$ws.Range("D51").Value = 203
$ws.Range("E51").Value = 2349
$ws.Range("F51").Value = 50
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 26

$ws.Range("A52").Value = "Tailandia"
$ws.Range("B52").Value = 2551
$ws.Range("C52").Value = 33
$ws.Range("D52").Value = 1218
$ws.Range("E52").Value = 1295
$ws.Range("F52").Value = 61
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 38

$ws.Range("A53").Value = "Singapur"
$ws.Range("B53").Value = 2299
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 528
$ws.Range("E53").Value = 1763
$ws.Range("F53").Value = 31
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 8

# Kazajistan (row 75)
$ws.Range("D75").Value = 86
$ws.Range("E75").Value = 801

# Letonia (row 83)
$ws.Range("B83").Value = 651
$ws.Range("C83").Value = 21
$ws.Range("E83").Value = 630
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 5

# Malta (row 103)
$ws.Range("D103").Value = 44
$ws.Range("E103").Value = 323

# Sri Lanka (row 116)
$ws.Range("D116").Value = 55
$ws.Range("E116").Value = 137
